$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3938.2307
$ws.Range("J40").Value = 4016.4167
$ws.Range("L40").Value = 4016.4167
$ws.Range("N40").Value = -4366.4167

$ws.Range("H43").Value = 5248.75
$ws.Range("J43").Value = 5427.143
$ws.Range("L43").Value = 5427.143
$ws.Range("N43").Value = -5565.143

$ws.Range("H112").Value = 2567487.2
$ws.Range("J112").Value = 2567487.2
$ws.Range("L112").Value = 7702461.600000001
$ws.Range("N112").Value = -7704677.600000001

$ws.Range("H132").Value = 1805.5577
$ws.Range("I132").Value = 1722.75
$ws.Range("K132").Value = 5168.25
$ws.Range("M132").Value = -2638.25

$ws.Range("H138").Value = 2387.889
$ws.Range("I138").Value = 916.8125
$ws.Range("K138").Value = 2750.4375
$ws.Range("M138").Value = 2389.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2744.2307
$ws.Range("I45").Value = 1537.2858
$ws.Range("K45").Value = 1537.2858
$ws.Range("M45").Value = -1160.2858

$ws.Range("H97").Value = 1873.3636
$ws.Range("I97").Value = 1873.3636
$ws.Range("K97").Value = 1873.3636
$ws.Range("M97").Value = -1377.3636

$ws.Range("H122").Value = 3064.9565
$ws.Range("I122").Value = 2150.0715
$ws.Range("K122").Value = 6450.2145
$ws.Range("M122").Value = -4000.2145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 104997.75
$ws.Range("J13").Value = 104997.75
$ws.Range("L13").Value = 104997.75
$ws.Range("N13").Value = -105333.75

$ws.Range("H99").Value = 3346.8823
$ws.Range("I99").Value = 2842.5715
$ws.Range("J99").Value = 3699.9
$ws.Range("K99").Value = 2842.5715
$ws.Range("L99").Value = 3699.9
$ws.Range("M99").Value = -1344.5715
$ws.Range("N99").Value = -6695.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5125.1665
$ws.Range("I31").Value = 2383
$ws.Range("J31").Value = 6179.846
$ws.Range("K31").Value = 2383
$ws.Range("L31").Value = 6179.846
$ws.Range("M31").Value = -2088
$ws.Range("N31").Value = -6769.846

$ws.Range("H34").Value = 5125.1665
$ws.Range("I34").Value = 2383
$ws.Range("J34").Value = 6179.846
$ws.Range("K34").Value = 2383
$ws.Range("L34").Value = 6179.846
$ws.Range("M34").Value = -2181
$ws.Range("N34").Value = -6583.846

$ws.Range("H122").Value = 3452887.8
$ws.Range("I122").Value = 5266400.5
$ws.Range("K122").Value = 15799201.5
$ws.Range("M122").Value = -15796751.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 2999
$ws.Range("I87").Value = 2999
$ws.Range("K87").Value = 8997
$ws.Range("M87").Value = -7749

$ws.Range("H90").Value = 2999
$ws.Range("I90").Value = 2999
$ws.Range("K90").Value = 26991
$ws.Range("M90").Value = -20751

$ws.Range("H131").Value = 1728.2
$ws.Range("J131").Value = 1747
$ws.Range("L131").Value = 5241
$ws.Range("N131").Value = -15321

$ws.Range("H132").Value = 836932.25
$ws.Range("I132").Value = 5500
$ws.Range("K132").Value = 49500
$ws.Range("M132").Value = -46970

$ws.Range("H134").Value = 5150.6665
$ws.Range("I134").Value = 5150.6665
$ws.Range("K134").Value = 15451.9995
$ws.Range("M134").Value = -10381.9995

$ws.Range("H140").Value = 4750
$ws.Range("I140").Value = 3000
$ws.Range("K140").Value = 9000
$ws.Range("M140").Value = -3820

$ws.Range("H141").Value = 3879.077
$ws.Range("I141").Value = 2553.5
$ws.Range("J141").Value = 6000
$ws.Range("K141").Value = 7660.5
$ws.Range("L141").Value = 18000
$ws.Range("M141").Value = -2480.5
$ws.Range("N141").Value = -28360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 27369.926
$ws.Range("I70").Value = 105581.336
$ws.Range("J70").Value = 5023.8096
$ws.Range("K70").Value = 105581.336
$ws.Range("L70").Value = 5023.8096
$ws.Range("M70").Value = -105311.336
$ws.Range("N70").Value = -5563.8096

$ws.Range("H73").Value = 27369.926
$ws.Range("I73").Value = 105581.336
$ws.Range("J73").Value = 5023.8096
$ws.Range("K73").Value = 105581.336
$ws.Range("L73").Value = 5023.8096
$ws.Range("M73").Value = -104645.336
$ws.Range("N73").Value = -6895.8096

$ws.Range("H122").Value = 1861.2222
$ws.Range("I122").Value = 1624
$ws.Range("J122").Value = 2051
$ws.Range("K122").Value = 4872
$ws.Range("L122").Value = 6153
$ws.Range("M122").Value = -2422
$ws.Range("N122").Value = -11053

$ws.Range("H132").Value = 2055.0625
$ws.Range("I132").Value = 1582.25
$ws.Range("J132").Value = 3473.5
$ws.Range("K132").Value = 4746.75
$ws.Range("L132").Value = 10420.5
$ws.Range("M132").Value = -2216.75
$ws.Range("N132").Value = -15480.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2728.1428
$ws.Range("I93").Value = 2439.8
$ws.Range("J93").Value = 3449
$ws.Range("K93").Value = 2439.8
$ws.Range("L93").Value = 3449
$ws.Range("M93").Value = -1191.8
$ws.Range("N93").Value = -5945

$ws.Range("H122").Value = 21470.2
$ws.Range("I122").Value = 22079.084
$ws.Range("J122").Value = 19034.666
$ws.Range("K122").Value = 66237.25199999999
$ws.Range("L122").Value = 57103.99800000001
$ws.Range("M122").Value = -63787.25199999999
$ws.Range("N122").Value = -62003.99800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 404999.5
$ws.Range("I5").Value = 9999
$ws.Range("K5").Value = 9999
$ws.Range("M5").Value = -9887

$ws.Range("H62").Value = 5197.4546
$ws.Range("I62").Value = 3849.3333
$ws.Range("J62").Value = 6130.769
$ws.Range("K62").Value = 3849.3333
$ws.Range("L62").Value = 6130.769
$ws.Range("M62").Value = -3225.3333
$ws.Range("N62").Value = -7378.769

$ws.Range("H65").Value = 5197.4546
$ws.Range("I65").Value = 3849.3333
$ws.Range("J65").Value = 6130.769
$ws.Range("K65").Value = 19246.6665
$ws.Range("L65").Value = 30653.845
$ws.Range("M65").Value = -16126.6665
$ws.Range("N65").Value = -36893.845

$ws.Range("H122").Value = 111118450
$ws.Range("I122").Value = 142865140
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 428595420
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -428592970
$ws.Range("N122").Value = -19900
